$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$targetRows = @(2,6,7,12,16,17,20,21,24,30,32,35,36,40,42,49,50,52,54,56,63,65,69,74,76,77,80,82,84,91,95,99,106,107,109,110,113,115,117,119,126,127,129,131,136,140,143,146,151,154,157,160,161,163,164,170,172,173,174,177,183,185,186,189)

foreach ($r in $targetRows) {
    $ws.Range("G" + $r + ":N" + $r).Value = "N/A"
}
